# Appends new daily price rows (2024-08-28 .. 2024-09-25) to the FMCG
# stocks worksheet, continuing the existing table starting at row 634.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("2024-08-28", 661.9000244140625,  1477.5,              611.2000122070312,  1308.050048828125,  839.75,              18182.55035400391, 0,                       170.8814800377791),
    @("2024-08-29", 660.75,             1476.699951171875,   603.6199951171875,  1309.849975585938,  813,                 18047.759765625,   -0.007413183835854168,   169.6147042121162),
    @("2024-08-30", 647.1500244140625,  1481.199951171875,   600.3599853515625,  1311.599975585938,  813.2000122070312,   17962.32995605469, -0.004733540931380745,   168.8118260671641),
    @("2024-09-02", 650.9500122070312,  1470.050048828125,   608.5800170898438,  1316.800048828125,  811.2000122070312,   17987.74047851562, 0.001414656256905703,    169.0506367731497),
    @("2024-09-03", 640.0499877929688,  1460.75,              599.9400024414062,  1341.949951171875,  822.3499755859375,  17977.66967773438, -0.0005598702512568748, 168.9559903506644),
    @("2024-09-04", 645.5999755859375,  1475.300048828125,   609,                1327.75,            824.2999877929688,   18052.54992675781, 0.004165181047695957,    169.6597226395676),
    @("2024-09-05", 643.9000244140625,  1457.699951171875,   602.1799926757812,  1254.800048828125,  835.4000244140625,   17792.94024658203, -0.01438077619112318,    167.21988413964),
    @("2024-09-06", 665.25,             1443.449951171875,   597.2999877929688,  1256.849975585938,  832.7000122070312,   17880.34979248047, 0.004912597057432855,    168.0413680504086),
    @("2024-09-09", 676,                1492.050048828125,   610.3400268554688,  1225.25,            827.5999755859375,   18025.32012939453, 0.008107802061849448,    169.4038142007637),
    @("2024-09-10", 680,                1503.050048828125,   608,                1246,               824.75,              18130.15014648438, 0.005815709032478913,    170.3890174931475),
    @("2024-09-11", 680.4500122070312,  1499.949951171875,   627.6599731445312,  1229,               814,                 18088.97985839844, -0.002270818926114677,   170.0020948874219),
    @("2024-09-12", 686.0999755859375,  1513.449951171875,   645.5999755859375,  1224.849975585938,  809.7000122070312,   18193.19958496094, 0.00576150382046627,     170.9815626066031),
    @("2024-09-13", 681.9500122070312,  1491.300048828125,   646.6500244140625,  1229.300048828125,  788.0499877929688,   18027.60040283203, -0.009102257211853799,   169.4252444452731),
    @("2024-09-16", 695.2000122070312,  1456.349975585938,   621.0499877929688,  1219.699951171875,  751.9500122070312,   17765.49987792969, -0.01453884704817227,    166.9619967301841),
    @("2024-09-17", 692,                1459.400024414062,   649.6500244140625,  1222.949951171875,  746.75,              17827,             0.003461772676980224,    167.5399812085587),
    @("2024-09-18", 695.2999877929688,  1432.150024414062,   646.7000122070312,  1224.550048828125,  744.5999755859375,   17755.70007324219, -0.003999547133999692,   166.8698971568857),
    @("2024-09-19", 697,                1444.849975585938,   652.1500244140625,  1197.849975585938,  747.2000122070312,   17752.34997558594, -0.0001886773060161447,  166.8384125942349),
    @("2024-09-20", 709,                1456.599975585938,   654.4500122070312,  1206.300048828125,  747.5499877929688,   17905.25006103516, 0.008612949027001824,    168.275383337655),
    @("2024-09-23", 702.5,              1449.300048828125,   654.0999755859375,  1190,               763.75,              17852.70007324219, -0.00293489270542646,    167.7815131425944),
    @("2024-09-24", 705.0999755859375,  1446.349975585938,   646.8499755859375,  1194.699951171875,  760.9500122070312,   17843.19958496094, -0.0005321597429113499, 167.6922265756952),
    @("2024-09-25", 689.2000122070312,  1429.550048828125,   633.2999877929688,  1175.349975585938,  742.5499877929688,   17509.20007324219, -0.0187185885652627,     164.553264780832)
)

$startRow = 634
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    $dateCell = $ws.Cells.Item($r, 1)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $data[0]       # A: Date (text)

    $ws.Cells.Item($r, 3).Value = $data[1]       # C
    $ws.Cells.Item($r, 4).Value = $data[2]       # D
    $ws.Cells.Item($r, 5).Value = $data[3]       # E
    $ws.Cells.Item($r, 6).Value = $data[4]       # F
    $ws.Cells.Item($r, 7).Value = $data[5]       # G
    $ws.Cells.Item($r, 8).Value = $data[6]       # H
    $ws.Cells.Item($r, 9).Value = $data[7]       # I
    $ws.Cells.Item($r, 10).Value = $data[8]      # J
}
